{"js": "// 1) \"Set cross plane field size as 200mm\" paragraph becomes a hyperlink line,\n//    followed by two new paragraphs (folder-location note, PDD note).\nconst setCrossResults = context.document.body.search(\"Set cross plane field size as 200mm\", { matchCase: true });\nsetCrossResults.load(\"text\");\nawait context.sync();\n\nif (setCrossResults.items.length === 0) {\n  throw new Error(\"Could not find 'Set cross plane field size as 200mm' paragraph\");\n}\n\nconst targetParagraph = setCrossResults.items[0].paragraphs.getFirst();\n\n// Replace the paragraph text with the hyperlink URL text, then turn the whole\n// paragraph range into a hyperlink, matching the Hyperlink run style used elsewhere.\ntargetParagraph.insertText(\"https://github.com/brianmanderson/EDWProfile\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst hyperlinkRange = targetParagraph.getRange();\nhyperlinkRange.hyperlink = \"https://github.com/brianmanderson/EDWProfile\";\nawait context.sync();\n\n// Trailing space after the hyperlink (its own, non-hyperlinked run).\ntargetParagraph.insertText(\" \", Word.InsertLocation.end);\nawait context.sync();\n\n// New paragraph: \"Set the folder location with ONLY ONE DOSE DICOM present (path variable).\"\nconst folderParagraph = targetParagraph.insertParagraph(\"Set the folder location with ONLY ONE DOSE DICOM present (path variable).\", Word.InsertLocation.after);\nawait context.sync();\n\n// New paragraph: \"Change PDD data to reflect current machine/clinic\"\n// (inserted before the bold/italic formatting below so it does not inherit it)\nconst pddParagraph = folderParagraph.insertParagraph(\"Change PDD data to reflect current machine/clinic\", Word.InsertLocation.after);\nawait context.sync();\n\nconst boldItalicResults = folderParagraph.search(\"ONLY ONE DOSE DICOM\", { matchCase: true });\nboldItalicResults.load(\"text\");\nawait context.sync();\nboldItalicResults.items[0].font.bold = true;\nboldItalicResults.items[0].font.italic = true;\nawait context.sync();\n\n// 2) \"Calculating Wedge Angle\" heading now starts on a new page.\nconst headingResults = context.document.body.search(\"Calculating Wedge Angle\", { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\nconst headingParagraph = headingResults.items[0].paragraphs.getFirst();\nheadingParagraph.paragraphFormat.pageBreakBefore = true;\nawait context.sync();\n\n// 3) 15MV (60 degree) paragraph gains an explicit \" 60 degrees\" before the colon,\n//    and a brand-new 15MV 30-degree paragraph is added right after it.\nconst mv15Results = context.document.body.search(\"15MV: From the equations listed above\", { matchCase: true });\nmv15Results.load(\"text\");\nawait context.sync();\nconst mv15Paragraph = mv15Results.items[0].paragraphs.getFirst();\nmv15Paragraph.load(\"text\");\nawait context.sync();\nmv15Paragraph.insertText(\n  mv15Paragraph.text.replace(\"15MV:\", \"15MV 60 degrees:\"),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst mv15_30Paragraph = mv15Paragraph.insertParagraph(\n  \"15MV 30 degrees: From the equations listed above, and the exported dose profile, we would expect the measured angle to be 32.72 degrees, and between 29.93-35.35 degrees. The IC Profiler measured a response of xxx degrees.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 4) 6MV (60 degree) paragraph gains an explicit \" 60 degrees\" before the colon,\n//    and a brand-new 6MV 30-degree paragraph is added right after it.\nconst mv6Results = context.document.body.search(\"6MV: From the equations listed above\", { matchCase: true });\nmv6Results.load(\"text\");\nawait context.sync();\nconst mv6Paragraph = mv6Results.items[0].paragraphs.getFirst();\nmv6Paragraph.load(\"text\");\nawait context.sync();\nmv6Paragraph.insertText(\n  mv6Paragraph.text.replace(\"6MV:\", \"6MV 60 degrees:\"),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst mv6_30Paragraph = mv6Paragraph.insertParagraph(\n  \"6MV 30 degrees: From the equations listed above, and the exported dose profile, we would expect the measured angle to be 33.29 degrees, and between 31.01-35.46 degrees. The IC Profiler measured a response of xxx degrees.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 5) One extra blank paragraph at the very end of the document.\nmv6_30Paragraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexExact($doc, $exactText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\") -eq $exactText) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Find-ParagraphIndexStartsWith($doc, $prefixText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($prefixText)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# 1) \"Set cross plane field size as 200mm\" becomes a hyperlink to the code repo,\n#    followed by two brand-new paragraphs.\n$crossIdx = Find-ParagraphIndexExact $d \"Set cross plane field size as 200mm\"\n$crossParagraph = $d.Paragraphs.Item($crossIdx)\n$crossParagraph.Range.Text = \"https://github.com/brianmanderson/EDWProfile\"\n\n$linkRange = $d.Range($crossParagraph.Range.Start, $crossParagraph.Range.End - 1)\n$d.Hyperlinks.Add($linkRange, \"https://github.com/brianmanderson/EDWProfile\") | Out-Null\n\n# trailing space after the hyperlink, still inside the same paragraph\n$crossParagraph = $d.Paragraphs.Item($crossIdx)\n$trailingSpace = $d.Range($crossParagraph.Range.End - 1, $crossParagraph.Range.End - 1)\n$trailingSpace.InsertAfter(\" \")\n\n# new paragraph: folder-location note, with \"ONLY ONE DOSE DICOM\" bold + italic\n$crossParagraph = $d.Paragraphs.Item($crossIdx)\n$crossParagraph.Range.InsertParagraphAfter()\n$folderParagraph = $d.Paragraphs.Item($crossIdx + 1)\n$folderParagraph.Range.Text = \"Set the folder location with ONLY ONE DOSE DICOM present (path variable).\"\n\n# new paragraph: PDD note (inserted before the bold/italic formatting below so it\n# does not inherit any character formatting)\n$folderParagraph.Range.InsertParagraphAfter()\n$pddParagraph = $d.Paragraphs.Item($crossIdx + 2)\n$pddParagraph.Range.Text = \"Change PDD data to reflect current machine/clinic\"\n\n$folderParagraph = $d.Paragraphs.Item($crossIdx + 1)\n$boldItalicRange = $folderParagraph.Range.Duplicate\n$boldItalicRange.Find.Execute(\"ONLY ONE DOSE DICOM\") | Out-Null\n$boldItalicRange.Bold = 1\n$boldItalicRange.Italic = 1\n\n# 2) \"Calculating Wedge Angle\" heading now starts on a new page.\n$headingIdx = Find-ParagraphIndexExact $d \"Calculating Wedge Angle\"\n$headingParagraph = $d.Paragraphs.Item($headingIdx)\n$headingParagraph.Format.PageBreakBefore = 1\n\n# 3) 15MV (60 degree) paragraph gains an explicit \" 60 degrees\" before the colon,\n#    and a brand-new 15MV 30-degree paragraph is added right after it.\n$mv15Idx = Find-ParagraphIndexStartsWith $d \"15MV:\"\n$mv15Paragraph = $d.Paragraphs.Item($mv15Idx)\n$mv15Find = $mv15Paragraph.Range.Duplicate\n$mv15Find.Find.Execute(\"15MV:\", $false, $false, $false, $false, $false, $true, 1, $false, \"15MV 60 degrees:\", 2) | Out-Null\n\n$mv15Paragraph = $d.Paragraphs.Item($mv15Idx)\n$mv15Paragraph.Range.InsertParagraphAfter()\n$mv15_30Paragraph = $d.Paragraphs.Item($mv15Idx + 1)\n$mv15_30Paragraph.Range.Text = \"15MV 30 degrees: From the equations listed above, and the exported dose profile, we would expect the measured angle to be 32.72 degrees, and between 29.93-35.35 degrees. The IC Profiler measured a response of xxx degrees.\"\n\n# 4) 6MV (60 degree) paragraph gains an explicit \" 60 degrees\" before the colon,\n#    and a brand-new 6MV 30-degree paragraph is added right after it.\n$mv6Idx = Find-ParagraphIndexStartsWith $d \"6MV:\"\n$mv6Paragraph = $d.Paragraphs.Item($mv6Idx)\n$mv6Find = $mv6Paragraph.Range.Duplicate\n$mv6Find.Find.Execute(\"6MV:\", $false, $false, $false, $false, $false, $true, 1, $false, \"6MV 60 degrees:\", 2) | Out-Null\n\n$mv6Paragraph = $d.Paragraphs.Item($mv6Idx)\n$mv6Paragraph.Range.InsertParagraphAfter()\n$mv6_30Paragraph = $d.Paragraphs.Item($mv6Idx + 1)\n$mv6_30Paragraph.Range.Text = \"6MV 30 degrees: From the equations listed above, and the exported dose profile, we would expect the measured angle to be 33.29 degrees, and between 31.01-35.46 degrees. The IC Profiler measured a response of xxx degrees.\"\n\n# 5) One extra blank paragraph at the very end of the document.\n$mv6_30Paragraph = $d.Paragraphs.Item($mv6Idx + 1)\n$mv6_30Paragraph.Range.InsertParagraphAfter()\n"}
